# Actualización desde MV -datos-
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
# (equivalently: $ws = $wb.ActiveSheet)

# New row 15: date label must stay a text value (not auto-converted to a date serial)
$ws.Range("A15").NumberFormat = "@"
$ws.Range("A15").Value = "02-11-2021"
$ws.Range("A15").Style = "Normal"

$ws.Range("B15").Value = 3.75
$ws.Range("C15").Value = 4.5
$ws.Range("D15").Value = 5
$ws.Range("E15").Value = 5.25
$ws.Range("F15").Value = 5.5
$ws.Range("G15").Value = 5.5
$ws.Range("H15").Value = 5
$ws.Range("M15").Value = 9
